# P4-2136 test data obfuscating exercise.
#
# This script:
#  1) Obfuscates the "Bedford" test-fixture location data (Courts + Prisons
#     sheets) to a fictitious "Fred"-themed location, matching the new
#     shared-string content introduced by the commit (including splitting
#     the old Courts address into a new standalone "England" string and a
#     new Prisons postcode "FD40 1HG").
#  2) Re-settles each sheet's lingering "2:22" legacy selection so it just
#     tracks the sheet's active cell (no real range is selected any more).
#  3) Moves the active tab from "Police" to "Prisons", with Prisons' own
#     selection landing on K2 and Courts' own landing on C2.
#
# Note: view-only cosmetics (per-pane scroll/topLeftCell) are not part of
# the workbook's logical state and are left to the host to manage.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Cell value obfuscation
# ---------------------------------------------------------------------

# Courts sheet - "Bedford County Court" row
$courts = $wb.Worksheets.Item("Courts")
$courts.Range("C2").Value = "Freds County Court"
$courts.Range("D2").Value = "FRDCT"
$courts.Range("F2").Value = "Freds County Court"
$courts.Range("J2").Value = "Fredford"
$courts.Range("K2").Value = "Fredfordshire"
$courts.Range("L2").Value = "England"
$courts.Range("M2").Value = "FR31 3ZZ"

# Prisons sheet - "HMP Bedford" row
$prisons = $wb.Worksheets.Item("Prisons")
$prisons.Range("C2").Value = "HMP Fred"
$prisons.Range("F2").Value = "HMP Fred"
$prisons.Range("H2").Value = "Fred"
$prisons.Range("I2").Value = "Fredfordshire"
$prisons.Range("J2").Value = "Fred St"
$prisons.Range("K2").Value = "FD40 1HG"

# ---------------------------------------------------------------------
# 2) Per-sheet selection clean-up (stray "2:22" sqref collapses onto the
#    sheet's own active cell) - every sheet except JPCNOMIS, which the
#    diff leaves untouched.
# ---------------------------------------------------------------------

$wb.Worksheets.Item("QUERIES").Range("D18").Select() | Out-Null
$wb.Worksheets.Item("JPCU").Range("A2").Select() | Out-Null
$wb.Worksheets.Item("NOMIS ").Range("A1").Select() | Out-Null
$wb.Worksheets.Item("Overview").Range("A1").Select() | Out-Null
$wb.Worksheets.Item("Police").Range("A2").Select() | Out-Null
$wb.Worksheets.Item("Police Info Sheet").Range("A2").Select() | Out-Null
$wb.Worksheets.Item("Hospitals").Range("H2").Select() | Out-Null
$wb.Worksheets.Item("Immigration").Range("A2").Select() | Out-Null
$wb.Worksheets.Item("STC&SCH").Range("A2").Select() | Out-Null
$wb.Worksheets.Item("Other").Range("A2").Select() | Out-Null
$wb.Worksheets.Item("Update Sheet").Range("C35").Select() | Out-Null

# Courts moves its active cell from A3 to C2.
$wb.Worksheets.Item("Courts").Range("C2").Select() | Out-Null

# Prisons becomes the active tab, selection lands on K2 (its postcode
# cell that was just edited above). Select this one last so it ends up
# driving the workbook's activeTab.
$wb.Worksheets.Item("Prisons").Range("K2").Select() | Out-Null
